$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# The export is a rolling window of dates. This update drops the two oldest
# dates (2025-09-13 and 2025-09-14) and appends the newest date (2025-12-12)
# with a zero count, keeping the table the same width (90 data rows).

# Drop the two oldest date rows - this shifts every remaining row up by two
# and naturally trims the used range from C91 to C90.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# Append the new date row at the bottom of the (now shorter) table.
$lastRow = $ws.UsedRange.Row() + $ws.UsedRange.Rows.Count()

# Build the date as text via a formula so Excel's automatic date detection
# (which would turn a literal "2025-12-12" value into a date serial number)
# never kicks in, then freeze the formula result back down to a plain value.
$ws.Cells.Item($lastRow, 1).Formula = '=""&"2025-12-12"'
$ws.Cells.Item($lastRow, 1).Copy()
$ws.Cells.Item($lastRow, 1).PasteSpecial("xlPasteValues")

$ws.Cells.Item($lastRow, 2).Value = 0
$ws.Cells.Item($lastRow, 3).Value = 0

Write-Output "updated GSC export rolling window"
